$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.900730666666667
$ws.Range("N2").Value = 17.702192
$ws.Range("O2").Value = 0.03970749001357476
$ws.Range("P2").Value = 0.03970749001357476
$ws.Range("Q2").Value = 4.408547994949333
$ws.Range("R2").Value = 39.676931954544
$ws.Range("S2").Value = 0.001229723039603797
$ws.Range("T2").Value = 0.001229723039603797

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 111.5917106666667
$ws.Range("N3").Value = 334.775132
$ws.Range("O3").Value = 0.7509284844884279
$ws.Range("P3").Value = 0.7509284844884279
$ws.Range("Q3").Value = 83.37228728156933
$ws.Range("R3").Value = 750.3505855341239
$ws.Range("S3").Value = 0.02325591615472267
$ws.Range("T3").Value = 0.02325591615472267

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2093640254979974
$ws.Range("P4").Value = 0.2093640254979974
$ws.Range("Q4").Value = 23.24476703282367
$ws.Range("R4").Value = 209.202903295413
$ws.Range("S4").Value = 0.006483909351386016
$ws.Range("T4").Value = 0.006483909351386015

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.900730666666667
$ws.Range("N5").Value = 17.702192
$ws.Range("O5").Value = 0.03970749001357476
$ws.Range("P5").Value = 0.03970749001357476
$ws.Range("Q5").Value = 116.5169626511983
$ws.Range("R5").Value = 1048.652663860784
$ws.Range("S5").Value = 0.03250131191516728
$ws.Range("T5").Value = 0.03250131191516728

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 111.5917106666667
$ws.Range("N6").Value = 334.775132
$ws.Range("O6").Value = 0.7509284844884279
$ws.Range("P6").Value = 0.7509284844884279
$ws.Range("Q6").Value = 2203.511381629685
$ws.Range("R6").Value = 19831.60243466717
$ws.Range("S6").Value = 0.61464879527763
$ws.Range("T6").Value = 0.61464879527763

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2093640254979974
$ws.Range("P7").Value = 0.2093640254979974
$ws.Range("Q7").Value = 614.3541264158216
$ws.Range("R7").Value = 5529.187137742394
$ws.Range("S7").Value = 0.1713683109710326
$ws.Range("T7").Value = 0.1713683109710326

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.900730666666667
$ws.Range("N8").Value = 17.702192
$ws.Range("O8").Value = 0.03970749001357476
$ws.Range("P8").Value = 0.03970749001357476
$ws.Range("Q8").Value = 21.425547149936
$ws.Range("R8").Value = 192.829924349424
$ws.Range("S8").Value = 0.005976455058803681
$ws.Range("T8").Value = 0.005976455058803681

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 111.5917106666667
$ws.Range("N9").Value = 334.775132
$ws.Range("O9").Value = 0.7509284844884279
$ws.Range("P9").Value = 0.7509284844884279
$ws.Range("Q9").Value = 405.1893898389559
$ws.Range("R9").Value = 3646.704508550604
$ws.Range("S9").Value = 0.1130237730560752
$ws.Range("T9").Value = 0.1130237730560752

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 31.11253633333333
$ws.Range("N10").Value = 93.337609
$ws.Range("O10").Value = 0.2093640254979974
$ws.Range("P10").Value = 0.2093640254979974
$ws.Range("Q10").Value = 112.969588313797
$ws.Range("R10").Value = 1016.726294824173
$ws.Range("S10").Value = 0.03151180517557882
$ws.Range("T10").Value = 0.03151180517557882
